$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MasterExecutor")

# Update the "Runmode" (column E) values for rows 2-7: "Yes" -> "no"
$ws.Range("E2:E7").Value = "no"

# Update the "Runmode" (column E) values for rows 13-21: "no" -> "Yes"
$ws.Range("E13:E21").Value = "Yes"

# Update the selection on the active sheet to match the new state
$ws.Range("E2:E8").Select()
